# Auto-generated Excel COM-interop script to apply numeric odds updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 7).Value = 2.48   # G3: 2.44 -> 2.48
$ws.Cells.Item(3, 23).Value = 1.67   # W3: 1.69 -> 1.67
$ws.Cells.Item(3, 31).Value = 1000   # AE3: 34 -> 1000
$ws.Cells.Item(4, 14).Value = 2.74   # N4: 1.66 -> 2.74
$ws.Cells.Item(4, 20).Value = 1.04   # T4: 1.01 -> 1.04
$ws.Cells.Item(4, 21).Value = 1.04   # U4: 1.01 -> 1.04
$ws.Cells.Item(5, 7).Value = 3.7   # G5: 3.75 -> 3.7
$ws.Cells.Item(5, 10).Value = 3.25   # J5: 3.2 -> 3.25
$ws.Cells.Item(5, 11).Value = 3.3   # K5: 3.25 -> 3.3
$ws.Cells.Item(5, 12).Value = 1.48   # L5: 0 -> 1.48
$ws.Cells.Item(5, 22).Value = 1.72   # V5: 0 -> 1.72
$ws.Cells.Item(5, 23).Value = 1.37   # W5: 0 -> 1.37
$ws.Cells.Item(5, 24).Value = 10   # X5: 10.5 -> 10
$ws.Cells.Item(5, 26).Value = 13.5   # Z5: 17.5 -> 13.5
$ws.Cells.Item(5, 27).Value = 32   # AA5: 36 -> 32
$ws.Cells.Item(5, 30).Value = 11.5   # AD5: 12 -> 11.5
$ws.Cells.Item(5, 31).Value = 30   # AE5: 32 -> 30
$ws.Cells.Item(5, 32).Value = 24   # AF5: 25 -> 24
$ws.Cells.Item(5, 33).Value = 16   # AG5: 17 -> 16
$ws.Cells.Item(5, 35).Value = 50   # AI5: 60 -> 50
$ws.Cells.Item(5, 36).Value = 75   # AJ5: 85 -> 75
$ws.Cells.Item(5, 37).Value = 50   # AK5: 55 -> 50
$ws.Cells.Item(5, 38).Value = 70   # AL5: 80 -> 70
$ws.Cells.Item(5, 39).Value = 140   # AM5: 150 -> 140
$ws.Cells.Item(5, 40).Value = 60   # AN5: 70 -> 60
$ws.Cells.Item(5, 41).Value = 27   # AO5: 28 -> 27
$ws.Cells.Item(6, 6).Value = 1.77   # F6: 1.78 -> 1.77
$ws.Cells.Item(6, 7).Value = 1.86   # G6: 1.87 -> 1.86
$ws.Cells.Item(6, 9).Value = 5.2   # I6: 5.3 -> 5.2
$ws.Cells.Item(6, 10).Value = 3.75   # J6: 3.8 -> 3.75
$ws.Cells.Item(6, 12).Value = 1.01   # L6: 0 -> 1.01
$ws.Cells.Item(6, 13).Value = 1.01   # M6: 0 -> 1.01
$ws.Cells.Item(6, 18).Value = 1.31   # R6: 0 -> 1.31
$ws.Cells.Item(6, 19).Value = 2.88   # S6: 0 -> 2.88
$ws.Cells.Item(6, 20).Value = 1.04   # T6: 0 -> 1.04
$ws.Cells.Item(6, 21).Value = 1.04   # U6: 0 -> 1.04
$ws.Cells.Item(6, 22).Value = 1.23   # V6: 0 -> 1.23
$ws.Cells.Item(6, 23).Value = 2.16   # W6: 0 -> 2.16
$ws.Cells.Item(6, 24).Value = 18.5   # X6: 0 -> 18.5
$ws.Cells.Item(6, 25).Value = 21   # Y6: 0 -> 21
$ws.Cells.Item(6, 26).Value = 46   # Z6: 0 -> 46
$ws.Cells.Item(6, 27).Value = 1000   # AA6: 0 -> 1000
$ws.Cells.Item(6, 28).Value = 10.5   # AB6: 0 -> 10.5
$ws.Cells.Item(6, 29).Value = 10.5   # AC6: 0 -> 10.5
$ws.Cells.Item(6, 30).Value = 29   # AD6: 0 -> 29
$ws.Cells.Item(6, 31).Value = 80   # AE6: 0 -> 80
$ws.Cells.Item(6, 32).Value = 13   # AF6: 0 -> 13
$ws.Cells.Item(6, 33).Value = 12   # AG6: 0 -> 12
$ws.Cells.Item(6, 34).Value = 23   # AH6: 0 -> 23
$ws.Cells.Item(6, 35).Value = 80   # AI6: 0 -> 80
$ws.Cells.Item(6, 36).Value = 23   # AJ6: 0 -> 23
$ws.Cells.Item(6, 37).Value = 28   # AK6: 0 -> 28
$ws.Cells.Item(6, 38).Value = 42   # AL6: 0 -> 42
$ws.Cells.Item(6, 39).Value = 1000   # AM6: 0 -> 1000
$ws.Cells.Item(6, 40).Value = 14   # AN6: 0 -> 14
$ws.Cells.Item(6, 41).Value = 1000   # AO6: 0 -> 1000
$ws.Cells.Item(7, 7).Value = 6   # G7: 6.2 -> 6
$ws.Cells.Item(7, 9).Value = 1.62   # I7: 1.63 -> 1.62
$ws.Cells.Item(7, 10).Value = 4.7   # J7: 4.6 -> 4.7
$ws.Cells.Item(7, 11).Value = 4.8   # K7: 4.7 -> 4.8
$ws.Cells.Item(7, 12).Value = 1.28   # L7: 0 -> 1.28
$ws.Cells.Item(7, 16).Value = 2.42   # P7: 2.44 -> 2.42
$ws.Cells.Item(7, 18).Value = 1.56   # R7: 1.55 -> 1.56
$ws.Cells.Item(7, 20).Value = 1.74   # T7: 1.73 -> 1.74
$ws.Cells.Item(7, 22).Value = 2.6   # V7: 0 -> 2.6
$ws.Cells.Item(7, 23).Value = 1.2   # W7: 0 -> 1.2
$ws.Cells.Item(7, 26).Value = 10.5   # Z7: 11 -> 10.5
$ws.Cells.Item(7, 27).Value = 15.5   # AA7: 16.5 -> 15.5
$ws.Cells.Item(7, 29).Value = 10   # AC7: 10.5 -> 10
$ws.Cells.Item(7, 30).Value = 9.4   # AD7: 9.6 -> 9.4
$ws.Cells.Item(7, 31).Value = 14.5   # AE7: 15 -> 14.5
$ws.Cells.Item(7, 32).Value = 50   # AF7: 55 -> 50
$ws.Cells.Item(7, 33).Value = 22   # AG7: 23 -> 22
$ws.Cells.Item(7, 35).Value = 27   # AI7: 28 -> 27
$ws.Cells.Item(7, 37).Value = 70   # AK7: 75 -> 70
$ws.Cells.Item(7, 41).Value = 7.2   # AO7: 7.4 -> 7.2
$ws.Cells.Item(8, 12).Value = 1.42   # L8: 0 -> 1.42
$ws.Cells.Item(8, 13).Value = 1.07   # M8: 0 -> 1.07
$ws.Cells.Item(8, 14).Value = 3.5   # N8: 0 -> 3.5
$ws.Cells.Item(8, 15).Value = 1.34   # O8: 0 -> 1.34
$ws.Cells.Item(8, 17).Value = 1.98   # Q8: 1.92 -> 1.98
$ws.Cells.Item(8, 18).Value = 1.34   # R8: 0 -> 1.34
$ws.Cells.Item(8, 19).Value = 3.55   # S8: 0 -> 3.55
$ws.Cells.Item(8, 20).Value = 1.63   # T8: 0 -> 1.63
$ws.Cells.Item(8, 21).Value = 1.93   # U8: 0 -> 1.93
$ws.Cells.Item(8, 22).Value = 1.62   # V8: 0 -> 1.62
$ws.Cells.Item(8, 23).Value = 1.46   # W8: 0 -> 1.46
$ws.Cells.Item(8, 24).Value = 15   # X8: 0 -> 15
$ws.Cells.Item(8, 25).Value = 11   # Y8: 0 -> 11
$ws.Cells.Item(8, 26).Value = 16.5   # Z8: 0 -> 16.5
$ws.Cells.Item(8, 27).Value = 38   # AA8: 0 -> 38
$ws.Cells.Item(8, 28).Value = 12   # AB8: 0 -> 12
$ws.Cells.Item(8, 29).Value = 7.8   # AC8: 0 -> 7.8
$ws.Cells.Item(8, 30).Value = 14.5   # AD8: 0 -> 14.5
$ws.Cells.Item(8, 31).Value = 29   # AE8: 0 -> 29
$ws.Cells.Item(8, 32).Value = 21   # AF8: 0 -> 21
$ws.Cells.Item(8, 33).Value = 13.5   # AG8: 0 -> 13.5
$ws.Cells.Item(8, 34).Value = 18   # AH8: 0 -> 18
$ws.Cells.Item(8, 35).Value = 44   # AI8: 0 -> 44
$ws.Cells.Item(8, 36).Value = 55   # AJ8: 0 -> 55
$ws.Cells.Item(8, 37).Value = 36   # AK8: 0 -> 36
$ws.Cells.Item(8, 38).Value = 46   # AL8: 0 -> 46
$ws.Cells.Item(8, 39).Value = 110   # AM8: 0 -> 110
$ws.Cells.Item(8, 40).Value = 34   # AN8: 0 -> 34
$ws.Cells.Item(8, 41).Value = 29   # AO8: 0 -> 29
$ws.Cells.Item(9, 6).Value = 1.17   # F9: 1.19 -> 1.17
$ws.Cells.Item(9, 7).Value = 1.23   # G9: 1.25 -> 1.23
$ws.Cells.Item(9, 8).Value = 1.1   # H9: 1.11 -> 1.1
$ws.Cells.Item(9, 9).Value = 980   # I9: 870 -> 980
$ws.Cells.Item(9, 10).Value = 1.2   # J9: 1.09 -> 1.2
$ws.Cells.Item(10, 6).Value = 1.09   # F10: 1.3 -> 1.09
$ws.Cells.Item(10, 7).Value = 1.44   # G10: 1.53 -> 1.44
$ws.Cells.Item(10, 8).Value = 3.15   # H10: 1.09 -> 3.15
$ws.Cells.Item(10, 9).Value = 980   # I10: 46 -> 980
$ws.Cells.Item(10, 10).Value = 5   # J10: 4.8 -> 5
$ws.Cells.Item(10, 16).Value = 2.32   # P10: 1.25 -> 2.32
$ws.Cells.Item(10, 17).Value = 1.42   # Q10: 1.56 -> 1.42
$ws.Cells.Item(12, 26).Value = 10   # Z12: 9.800000000000001 -> 10
$ws.Cells.Item(13, 15).Value = 1.43   # O13: 1.42 -> 1.43
$ws.Cells.Item(13, 16).Value = 1.75   # P13: 1.74 -> 1.75
$ws.Cells.Item(13, 17).Value = 2.26   # Q13: 2.28 -> 2.26
$ws.Cells.Item(13, 24).Value = 12   # X13: 11.5 -> 12
$ws.Cells.Item(13, 27).Value = 130   # AA13: 140 -> 130
$ws.Cells.Item(18, 7).Value = 2.06   # G18: 2.1 -> 2.06
$ws.Cells.Item(18, 11).Value = 3.6   # K18: 3.55 -> 3.6
$ws.Cells.Item(19, 6).Value = 2.62   # F19: 2.64 -> 2.62
$ws.Cells.Item(19, 8).Value = 2.86   # H19: 2.92 -> 2.86
$ws.Cells.Item(20, 6).Value = 1.7   # F20: 1.78 -> 1.7
$ws.Cells.Item(20, 7).Value = 2.12   # G20: 2.26 -> 2.12
$ws.Cells.Item(20, 8).Value = 3.5   # H20: 3.25 -> 3.5
$ws.Cells.Item(20, 9).Value = 8.4   # I20: 7 -> 8.4
$ws.Cells.Item(20, 10).Value = 3.45   # J20: 3.4 -> 3.45
$ws.Cells.Item(20, 11).Value = 8   # K20: 7.8 -> 8
$ws.Cells.Item(20, 16).Value = 2.16   # P20: 2.14 -> 2.16
$ws.Cells.Item(20, 17).Value = 1.47   # Q20: 1.55 -> 1.47
$ws.Cells.Item(21, 7).Value = 1.18   # G21: 1000 -> 1.18
$ws.Cells.Item(21, 8).Value = 16   # H21: 1.04 -> 16
$ws.Cells.Item(22, 8).Value = 2.4   # H22: 2.42 -> 2.4
$ws.Cells.Item(23, 6).Value = 8.4   # F23: 8.199999999999999 -> 8.4
$ws.Cells.Item(23, 7).Value = 9.6   # G23: 9.4 -> 9.6
$ws.Cells.Item(23, 9).Value = 1.43   # I23: 1.44 -> 1.43
$ws.Cells.Item(23, 10).Value = 5.2   # J23: 5.1 -> 5.2
$ws.Cells.Item(23, 14).Value = 5.4   # N23: 5.3 -> 5.4
$ws.Cells.Item(23, 15).Value = 1.19   # O23: 1.2 -> 1.19
$ws.Cells.Item(23, 16).Value = 2.5   # P23: 2.48 -> 2.5
$ws.Cells.Item(23, 17).Value = 1.62   # Q23: 1.64 -> 1.62
$ws.Cells.Item(23, 18).Value = 1.6   # R23: 1.59 -> 1.6
$ws.Cells.Item(23, 19).Value = 2.4   # S23: 2.48 -> 2.4
$ws.Cells.Item(23, 20).Value = 1.85   # T23: 1.84 -> 1.85
$ws.Cells.Item(23, 21).Value = 1.98   # U23: 2 -> 1.98
$ws.Cells.Item(23, 24).Value = 32   # X23: 30 -> 32
$ws.Cells.Item(23, 25).Value = 11   # Y23: 10.5 -> 11
$ws.Cells.Item(23, 26).Value = 10   # Z23: 9.800000000000001 -> 10
$ws.Cells.Item(23, 29).Value = 13.5   # AC23: 13 -> 13.5
$ws.Cells.Item(23, 35).Value = 34   # AI23: 1000 -> 34
$ws.Cells.Item(23, 36).Value = 330   # AJ23: 290 -> 330
$ws.Cells.Item(23, 37).Value = 150   # AK23: 140 -> 150
$ws.Cells.Item(23, 38).Value = 120   # AL23: 110 -> 120
$ws.Cells.Item(23, 39).Value = 140   # AM23: 130 -> 140
$ws.Cells.Item(23, 41).Value = 5.3   # AO23: 5.6 -> 5.3
$ws.Cells.Item(24, 6).Value = 1.82   # F24: 1.81 -> 1.82
$ws.Cells.Item(24, 9).Value = 5.2   # I24: 5.3 -> 5.2
$ws.Cells.Item(24, 10).Value = 3.7   # J24: 3.65 -> 3.7
$ws.Cells.Item(24, 17).Value = 1.94   # Q24: 1.97 -> 1.94
$ws.Cells.Item(25, 11).Value = 5.4   # K25: 5.3 -> 5.4
$ws.Cells.Item(25, 26).Value = 1000   # Z25: 80 -> 1000
$ws.Cells.Item(25, 28).Value = 9   # AB25: 9.199999999999999 -> 9
$ws.Cells.Item(25, 33).Value = 10.5   # AG25: 10 -> 10.5
$ws.Cells.Item(25, 37).Value = 14   # AK25: 14.5 -> 14
$ws.Cells.Item(26, 33).Value = 10   # AG26: 10.5 -> 10
$ws.Cells.Item(27, 33).Value = 13   # AG27: 12.5 -> 13
$ws.Cells.Item(27, 34).Value = 15   # AH27: 16 -> 15
$ws.Cells.Item(28, 6).Value = 1.26   # F28: 1.25 -> 1.26
$ws.Cells.Item(28, 7).Value = 1.27   # G28: 1.26 -> 1.27
$ws.Cells.Item(28, 8).Value = 13.5   # H28: 13 -> 13.5
$ws.Cells.Item(28, 9).Value = 14   # I28: 13.5 -> 14
$ws.Cells.Item(28, 19).Value = 2.38   # S28: 2.4 -> 2.38
$ws.Cells.Item(28, 29).Value = 17.5   # AC28: 17 -> 17.5
$ws.Cells.Item(28, 32).Value = 8   # AF28: 7.8 -> 8
$ws.Cells.Item(28, 36).Value = 9.199999999999999   # AJ28: 9 -> 9.199999999999999
$ws.Cells.Item(28, 38).Value = 36   # AL28: 38 -> 36
$ws.Cells.Item(28, 40).Value = 4   # AN28: 3.95 -> 4
$ws.Cells.Item(29, 18).Value = 1.51   # R29: 1.52 -> 1.51
$ws.Cells.Item(30, 6).Value = 4.3   # F30: 4.2 -> 4.3
$ws.Cells.Item(30, 7).Value = 4.4   # G30: 4.3 -> 4.4
$ws.Cells.Item(30, 10).Value = 3.45   # J30: 3.5 -> 3.45
$ws.Cells.Item(30, 11).Value = 3.5   # K30: 3.55 -> 3.5
$ws.Cells.Item(34, 6).Value = 3.3   # F34: 3.25 -> 3.3
$ws.Cells.Item(34, 7).Value = 3.4   # G34: 3.35 -> 3.4
$ws.Cells.Item(34, 8).Value = 2.54   # H34: 2.56 -> 2.54
$ws.Cells.Item(34, 9).Value = 2.58   # I34: 2.6 -> 2.58
$ws.Cells.Item(34, 14).Value = 3   # N34: 3.05 -> 3
$ws.Cells.Item(34, 21).Value = 1.93   # U34: 1.92 -> 1.93
$ws.Cells.Item(34, 25).Value = 8.6   # Y34: 8.800000000000001 -> 8.6
$ws.Cells.Item(34, 27).Value = 40   # AA34: 42 -> 40
$ws.Cells.Item(34, 28).Value = 10   # AB34: 9.800000000000001 -> 10
$ws.Cells.Item(34, 31).Value = 34   # AE34: 38 -> 34
$ws.Cells.Item(34, 33).Value = 15.5   # AG34: 14.5 -> 15.5
$ws.Cells.Item(34, 36).Value = 70   # AJ34: 75 -> 70
$ws.Cells.Item(34, 37).Value = 48   # AK34: 50 -> 48
$ws.Cells.Item(34, 38).Value = 70   # AL34: 75 -> 70
$ws.Cells.Item(34, 41).Value = 34   # AO34: 42 -> 34
$ws.Cells.Item(35, 9).Value = 4.4   # I35: 4.3 -> 4.4
$ws.Cells.Item(35, 24).Value = 19   # X35: 18.5 -> 19
$ws.Cells.Item(35, 37).Value = 17.5   # AK35: 16 -> 17.5
$ws.Cells.Item(36, 6).Value = 1.49   # F36: 1.52 -> 1.49
$ws.Cells.Item(36, 7).Value = 1.8   # G36: 1.81 -> 1.8
$ws.Cells.Item(37, 6).Value = 1.92   # F37: 1.84 -> 1.92
$ws.Cells.Item(37, 8).Value = 4.3   # H37: 4.6 -> 4.3
$ws.Cells.Item(37, 11).Value = 4.1   # K37: 4.2 -> 4.1
$ws.Cells.Item(38, 6).Value = 1.04   # F38: 1.09 -> 1.04
$ws.Cells.Item(38, 14).Value = 1.1   # N38: 1.08 -> 1.1
